# Updates the cryptos price/volume table (Sheet1, columns D "Price" and
# E "Volume(1h)") with refreshed values, matching the GitHub Actions
# "Updated cryptos list" commit.
#
# Column D values are textual (e.g. "62.633.12", using '.' as a thousands
# separator) so they must stay text cells, not get auto-parsed as numbers.
# Numeric-looking ones (e.g. "565.60") are written with a leading
# apostrophe to force text entry (as typing them in the Excel UI would),
# then ClearFormats() drops the resulting "Number Stored as Text"
# (quotePrefix) cell style since these cells carried no style originally.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.633.12"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "2.538.46"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'565.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'145.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.580"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").Value = "2.537.76"
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "'5.60"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'0.351"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "'26.93"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "2.994.10"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "62.650.75"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").Value = "2.542.15"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("D19").Value = "'11.39"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "'333.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").Value = "'6.72"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'64.46"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("E25").Value = "  -4.33%  "
$ws.Range("E26").Value = "  +3.57%  "
$ws.Range("E27").Value = "  +10.94%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "'8.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'7.21"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.92%  "
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "'176.32"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").Value = "'404.49"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.35%  "
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").Value = "'18.82"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("D40").Value = "'1.74"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'38.95"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.75%  "
$ws.Range("D43").Value = "'152.79"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "'3.72"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "'20.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").Value = "'0.0234"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").Value = "'18.15"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  -2.11%  "
